$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New shared string used by row 76 (column A)
# "Key appears" becomes sharedStrings index 22

# Data for new rows 76..91: ColumnA text, B, C, D(formula)
$newRows = @(
    @{ Row=76; A="Key appears"; B=120939; C=118643 },
    @{ Row=77; A="Get key";     B=121084; C=118788 },
    @{ Row=78; A="Room app";    B=121665; C=119246 },
    @{ Row=79; A="Get treasure";B=122084; C=119665 },
    @{ Row=80; A="Room app";    B=122678; C=120135 },
    @{ Row=81; A="Enter door";  B=124219; C=121676 },
    @{ Row=82; A="Room app";    B=124750; C=122084 },
    @{ Row=83; A="Room app";    B=125026; C=122360 },
    @{ Row=84; A="Get Key";     B=126650; C=123985 },
    @{ Row=85; A="Last heart";  B=$null;  C=125324 },
    @{ Row=86; A="Get key";     B=129103; C=126438 },
    @{ Row=87; A="Last heart";  B=129917; C=127252 },
    @{ Row=88; A="Get key";     B=130061; C=127396 },
    @{ Row=89; A="Get key";     B=130939; C=128274 },
    @{ Row=90; A="Get key";     B=131427; C=128763 },
    @{ Row=91; A="Get key";     B=133230; C=130565 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    if ($null -ne $r.B) {
        $ws.Cells.Item($row, 2).Value = $r.B
    }
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Formula = "=C$row-B$row"
}

# Update the view to match the post-edit state (scroll + selection)
$ws.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 70
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Window scroll position isn't always controllable via this host; ignore.
}
$ws.Range("C92").Select()
